$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.489.36'
$ws.Range('E2').Value = '  +0.02%  '

$ws.Range('D3').Value = '1.626.93'
$ws.Range('E3').Value = '  +0.13%  '

$ws.Range('E4').Value = '  +0.23%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.25%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.497'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.14%  '

$ws.Range('E7').Value = '  +0.20%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.250'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.41%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0620'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.45%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.84%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0842'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.72%  '

$ws.Range('D12').Value = '1.853.95'
$ws.Range('E12').Value = '  +0.16%  '

$ws.Range('D13').Value = '1.662.41'
$ws.Range('E13').Value = '  +2.18%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.64%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.16%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.46'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.10%  '

$ws.Range('D17').Value = '26.481.21'
$ws.Range('E17').Value = '  +0.05%  '

$ws.Range('D18').Value = '0.0₃0739'
$ws.Range('E18').Value = '  +1.24%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '214.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.89%  '

$ws.Range('E20').Value = '  +0.19%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.32%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.68%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.83%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.80%  '

$ws.Range('E26').Value = '  +0.16%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.120'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.78%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.80%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.62%  '

$ws.Range('E30').Value = '  -3.76%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.45%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.30'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.87%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.94'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.31%  '

$ws.Range('E34').Value = '  -0.21%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.38'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.27%  '

$ws.Range('D36').Value = '1.215.73'
$ws.Range('E36').Value = '  +5.39%  '

$ws.Range('E37').Value = '  +5.07%  '

$ws.Range('B38').Value = 'PaxDollar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.15%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.797'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.86%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.498'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.26%  '

$ws.Range('E41').Value = '  -1.63%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.794'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.31%  '

$ws.Range('E43').Value = '  -1.04%  '

$ws.Range('D44').Value = '1.762.50'
$ws.Range('E44').Value = '  +0.02%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.03%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.56'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.12%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.93'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.67%  '

$ws.Range('D48').Value = '0.0₆0102'
$ws.Range('E48').Value = '  -0.03%  '

$ws.Range('E49').Value = '  +0.44%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.61'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.41%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.408'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.32%  '
